# Apply the committed changes to MAI_holdings.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet is protected; unprotect so the values can be updated
$ws.Unprotect("D382")

# Update the confidential disclaimer date from 2021-03-23 to 2021-03-24
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-24 for illustrative purposes only and are subject to change."

# Update the Weight (D) and Percent Change (E) columns for rows 2-7
$ws.Range("D2").Value = 0.4933444174184345
$ws.Range("E2").Value = 0.003162055335968361

$ws.Range("D3").Value = 0.3292127220665661
$ws.Range("E3").Value = 0.0009118541033432681

$ws.Range("D4").Value = 0.09373492615741454
$ws.Range("E4").Value = -0.00249089863958607

$ws.Range("D5").Value = 0.05545891191231257
$ws.Range("E5").Value = 0.0005783021050198567

$ws.Range("D6").Value = 0.02824902244527233
$ws.Range("E6").Value = 0.006811989100817373

$ws.Range("E7").Value = 0.001851196157517165

# Re-protect the sheet with the original password and settings
$ws.Protect("D382", $false, $true, $true, $true)

